$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 46 (shifts the existing rows 46-73 down to 47-74)
$ws.Rows.Item(46).Insert()

# Populate the new row 46 with this week's data; most descriptive columns
# (A, B, C, E-L) repeat the values of the row directly below (same
# mercado/producto/variedad/calidad context), only the date + pricing
# columns change for the new entry.
$ws.Range("A46").Value = 1
$ws.Range("B46").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C46").Value = "Arica y Parinacota"
$ws.Range("D46").Value = 44586
$ws.Range("E46").Value = 15
$ws.Range("F46").Value = "Fruta"
$ws.Range("G46").Value = 100109
$ws.Range("H46").Value = "Uva"
$ws.Range("I46").Value = 100109001
$ws.Range("J46").Value = "Uva"
$ws.Range("K46").Value = "Superior Seedless"
$ws.Range("L46").Value = "Segunda"
$ws.Range("M46").Value = 250
$ws.Range("N46").Value = 19000
$ws.Range("O46").Value = 20000
$ws.Range("P46").Value = 19500
$ws.Range("Q46").Value = "$/caja 25 kilos"
$ws.Range("R46").Value = "Región de Coquimbo"
$ws.Range("S46").Value = 780
$ws.Range("T46").Value = 25
